# Video Script.docx edit:
#   1. "...gradually increase until the game becomes impossible" ->
#      "...gradually decrease until the game becomes impossible"
#   2. A new closing line "Thank you for watching" is added as its own
#      paragraph at the end of the document (reusing the slot of the old
#      trailing empty paragraph), with the "_GoBack" bookmark now wrapping
#      that new run instead of sitting at the end of the previous paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Word fix: increase -> decrease
# ---------------------------------------------------------------------------
$range = $d.Content
$range.Find.Execute(
    "interval between the boss’s shots gradually increase",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "interval between the boss’s shots gradually decrease", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Drop the old trailing empty paragraph by merging its paragraph mark
#    into the end of the "...impossible." paragraph, then re-insert a fresh
#    paragraph break in the same spot. This gives us a brand-new, bare
#    <w:p> (no leftover rsid/paraId bookkeeping) to hold the new sentence,
#    matching how the real edit was produced.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$boundary = $p3.Range.End - 1
$mergeRange = $d.Range($boundary, $boundary + 2)
$mergeRange.Delete()

$p3 = $d.Paragraphs.Item(3)
$splitPos = $p3.Range.End - 1
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark off the end of paragraph 3 and onto the
#    new "Thank you for watching" sentence in the new last paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$p4 = $d.Paragraphs.Item(4)
$insertion = $d.Range($p4.Range.Start, $p4.Range.Start)
$insertion.InsertBefore("Thank you for watching")

$p4 = $d.Paragraphs.Item(4)
$newBookmarkRange = $d.Range($p4.Range.Start, $p4.Range.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
